# Data cleaning updates:
# 1. Rename header "Copies Sold" -> "Copies_Sold"
# 2. Normalize "State" column values "Uttar pradesh" / "Uttar-Pradesh" -> "Uttar Pradesh"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header in F1
$headerCell = $ws.Range("F1")
if ($headerCell.Value() -eq "Copies Sold") {
    $headerCell.Value = "Copies_Sold"
}

# Determine the used range extent
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

# Column D holds the "State" values; normalize inconsistent spellings/hyphenation
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)
    $val = $cell.Value()
    if ($val -eq "Uttar pradesh" -or $val -eq "Uttar-Pradesh") {
        $cell.Value = "Uttar Pradesh"
    }
}
